$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.508.70"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "3.269.45"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.426"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("D12").Value = "3.836.71"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.138"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "68.537.38"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "3.277.86"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.517"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000120"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.42%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.830"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0689"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "344.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("D46").Value = "2.612.86"
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0283"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("E51").Value = "  -0.04%  "
